# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets
# to reflect newly generated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1096   # was 1094
$ws1.Range("F4").Value = 1716   # was 1713
$ws1.Range("F5").Value = 765    # was 764

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1096   # was 1094
$ws4.Range("F4").Value = 1716   # was 1713
$ws4.Range("F6").Value = 765    # was 764
